{"js": "// Replace the date line and every \"A\u00f7B=C, D\" answer cell with its new value.\n// Each old value is unique in the document, so an exact (case-sensitive,\n// non-wildcard) search-and-replace keyed on the old text is unambiguous.\nconst replacements = [\n  [\"2025-08-07 Thursday\", \"2025-08-08 Friday\"],\n  [\"714\u00f74=178, 2\", \"949\u00f72=474, 1\"],\n  [\"103\u00f78=12, 7\", \"386\u00f75=77, 1\"],\n  [\"570\u00f75=114, 0\", \"114\u00f74=28, 2\"],\n  [\"877\u00f79=97, 4\", \"554\u00f78=69, 2\"],\n  [\"369\u00f72=184, 1\", \"886\u00f79=98, 4\"],\n  [\"111\u00f72=55, 1\", \"465\u00f79=51, 6\"],\n  [\"602\u00f79=66, 8\", \"668\u00f78=83, 4\"],\n  [\"665\u00f75=133, 0\", \"260\u00f79=28, 8\"],\n  [\"944\u00f77=134, 6\", \"355\u00f73=118, 1\"],\n  [\"673\u00f76=112, 1\", \"752\u00f72=376, 0\"],\n  [\"190\u00f79=21, 1\", \"779\u00f73=259, 2\"],\n  [\"275\u00f79=30, 5\", \"257\u00f78=32, 1\"],\n  [\"613\u00f77=87, 4\", \"684\u00f78=85, 4\"],\n  [\"790\u00f77=112, 6\", \"165\u00f77=23, 4\"],\n  [\"491\u00f76=81, 5\", \"214\u00f72=107, 0\"],\n  [\"943\u00f77=134, 5\", \"831\u00f74=207, 3\"],\n  [\"754\u00f77=107, 5\", \"936\u00f79=104, 0\"],\n  [\"450\u00f79=50, 0\", \"357\u00f74=89, 1\"],\n  [\"510\u00f73=170, 0\", \"453\u00f74=113, 1\"],\n  [\"308\u00f79=34, 2\", \"155\u00f72=77, 1\"],\n  [\"203\u00f77=29, 0\", \"461\u00f72=230, 1\"],\n  [\"715\u00f74=178, 3\", \"224\u00f78=28, 0\"],\n  [\"898\u00f73=299, 1\", \"842\u00f73=280, 2\"],\n  [\"288\u00f77=41, 1\", \"408\u00f76=68, 0\"],\n  [\"435\u00f75=87, 0\", \"184\u00f77=26, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00f7B=C, D\" answer cell with its new value.\n# Each old value is unique in the document, so a case-sensitive Find/Replace\n# keyed on the old text unambiguously targets the right run.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-08-07 Thursday\"; New = \"2025-08-08 Friday\" }\n    @{ Old = \"714\u00f74=178, 2\"; New = \"949\u00f72=474, 1\" }\n    @{ Old = \"103\u00f78=12, 7\"; New = \"386\u00f75=77, 1\" }\n    @{ Old = \"570\u00f75=114, 0\"; New = \"114\u00f74=28, 2\" }\n    @{ Old = \"877\u00f79=97, 4\"; New = \"554\u00f78=69, 2\" }\n    @{ Old = \"369\u00f72=184, 1\"; New = \"886\u00f79=98, 4\" }\n    @{ Old = \"111\u00f72=55, 1\"; New = \"465\u00f79=51, 6\" }\n    @{ Old = \"602\u00f79=66, 8\"; New = \"668\u00f78=83, 4\" }\n    @{ Old = \"665\u00f75=133, 0\"; New = \"260\u00f79=28, 8\" }\n    @{ Old = \"944\u00f77=134, 6\"; New = \"355\u00f73=118, 1\" }\n    @{ Old = \"673\u00f76=112, 1\"; New = \"752\u00f72=376, 0\" }\n    @{ Old = \"190\u00f79=21, 1\"; New = \"779\u00f73=259, 2\" }\n    @{ Old = \"275\u00f79=30, 5\"; New = \"257\u00f78=32, 1\" }\n    @{ Old = \"613\u00f77=87, 4\"; New = \"684\u00f78=85, 4\" }\n    @{ Old = \"790\u00f77=112, 6\"; New = \"165\u00f77=23, 4\" }\n    @{ Old = \"491\u00f76=81, 5\"; New = \"214\u00f72=107, 0\" }\n    @{ Old = \"943\u00f77=134, 5\"; New = \"831\u00f74=207, 3\" }\n    @{ Old = \"754\u00f77=107, 5\"; New = \"936\u00f79=104, 0\" }\n    @{ Old = \"450\u00f79=50, 0\"; New = \"357\u00f74=89, 1\" }\n    @{ Old = \"510\u00f73=170, 0\"; New = \"453\u00f74=113, 1\" }\n    @{ Old = \"308\u00f79=34, 2\"; New = \"155\u00f72=77, 1\" }\n    @{ Old = \"203\u00f77=29, 0\"; New = \"461\u00f72=230, 1\" }\n    @{ Old = \"715\u00f74=178, 3\"; New = \"224\u00f78=28, 0\" }\n    @{ Old = \"898\u00f73=299, 1\"; New = \"842\u00f73=280, 2\" }\n    @{ Old = \"288\u00f77=41, 1\"; New = \"408\u00f76=68, 0\" }\n    @{ Old = \"435\u00f75=87, 0\"; New = \"184\u00f77=26, 2\" }\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $rng.Find.Execute($pair.Old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
